$d = $word.ActiveDocument

# 1. Title change (appears twice: main heading + bold run near end)
$d.Content.Find.Execute("Play Golden Castle Slot for Free - Review and Rating", $true, $false, $false, $false, $false, $true, 1, $false, "Play Golden Castle Slot for Free", 2)

# 2. "What we like" bullet items
$d.Content.Find.Execute("Innovative Titanways engine and exciting gameplay mechanics", $true, $false, $false, $false, $false, $true, 1, $false, "Innovative Titanways engine and unique symbol designs", 2)

$d.Content.Find.Execute("Avalanche feature adds extra excitement and chance for increased wins", $true, $false, $false, $false, $false, $true, 1, $false, "Avalanche feature adds excitement and potential for big wins", 2)

$d.Content.Find.Execute("Bonus game with progressive win multiplier and the potential for significant payouts", $true, $false, $false, $false, $false, $true, 1, $false, "Phantom Free Spins with progressive win multiplier", 2)

$d.Content.Find.Execute("Ghost King symbol and randomly transforming grid symbols offer unique gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Ghost King symbol adds an element of surprise to gameplay", 2)

# 3. "What we don't like" bullet items
$d.Content.Find.Execute("Highly volatile gameplay may not be suitable for all players", $true, $false, $false, $false, $false, $true, 1, $false, "Highly volatile gameplay may not suit all players", 2)

$d.Content.Find.Execute("Purchase of buy feature is expensive at 70 times default bet", $true, $false, $false, $false, $false, $true, 1, $false, "Buy feature can be costly", 2)

# 4. Meta description paragraph
$d.Content.Find.Execute("Read our neutral Golden Castle Slot review to play this innovative game with Titanways engine and Ghost King feature for free. Highly volatile with a top prize of 20,000x.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Golden Castle, an innovative slot game with exciting features. Play for free now!", 2)
